$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 27, shifting existing rows 27..48 down to 28..49.
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new weekly record.
$ws.Cells.Item(27, 1).Value = 10
$ws.Cells.Item(27, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(27, 3).Value = "La Araucanía"
$ws.Cells.Item(27, 4).Value = 44494
$ws.Cells.Item(27, 4).NumberFormat = $ws.Cells.Item(28, 4).NumberFormat
$ws.Cells.Item(27, 5).Value = 9
$ws.Cells.Item(27, 6).Value = 100112022
$ws.Cells.Item(27, 7).Value = "Arveja Verde"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 30
$ws.Cells.Item(27, 11).Value = 21000
$ws.Cells.Item(27, 12).Value = 21000
$ws.Cells.Item(27, 13).Value = 21000
$ws.Cells.Item(27, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(27, 15).Value = "Región Metropolitana"
$ws.Cells.Item(27, 16).Value = 840
$ws.Cells.Item(27, 17).Value = 25
$ws.Cells.Item(27, 18).Value = "Hortaliza"
